# NIT-9015338857.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos
# estado de cuenta"
#
# The account-statement table (rows 16-29 of Hoja1) previously listed all
# mora-periods for worker 1 (22802542 / KARETH LAVINIA GUZMAN PAJARO) first,
# then all periods for worker 2 (1091680446 / KELLY DAYANA ASCANIO TORRES).
# The update interleaves the two workers period-by-period (worker1/worker2
# pairs per periodo, ascending 2305..2311), refreshes the "Valor Mora" (F)
# figures and drops "Salario Basico" (G) from 1,160,000 to 1,000,000 for
# every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2305"; F = 46400; G = 1000000 },
    @{ Row = 17; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2305"; F = 46400; G = 1000000 },
    @{ Row = 18; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2306"; F = 46400; G = 1000000 },
    @{ Row = 19; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2306"; F = 46400; G = 1000000 },
    @{ Row = 20; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2307"; F = 46400; G = 1000000 },
    @{ Row = 21; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2307"; F = 46400; G = 1000000 },
    @{ Row = 22; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2308"; F = 46400; G = 1000000 },
    @{ Row = 23; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2308"; F = 46400; G = 1000000 },
    @{ Row = 24; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2309"; F = 46400; G = 1000000 },
    @{ Row = 25; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2309"; F = 46400; G = 1000000 },
    @{ Row = 26; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2310"; F = 46400; G = 1000000 },
    @{ Row = 27; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2310"; F = 46400; G = 1000000 },
    @{ Row = 28; C = "22802542";   D = "KARETH LAVINIA GUZMAN PAJARO"; E = "2311"; F = 37333; G = 1000000 },
    @{ Row = 29; C = "1091680446"; D = "KELLY DAYANA ASCANIO TORRES";  E = "2311"; F = 37333; G = 1000000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $item.D   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $item.E   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $item.F   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $item.G   # G - Salario Basico
}
